# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right before the "总计" sheet and
#    populate it with the same column layout used by the other quarterly
#    sheets (2020-Q4 / 2021-Q1).
# 2) Update the "总计" (totals) sheet: add a new top data row for 2022-Q1
#    and shift the existing rows down, renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted before "总计"
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totals)
$q1.Name = "2022-Q1"

# Headers (row 1)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data (row 2) - keep the textual columns as text so things like fund
# codes retain formatting (mirrors the other quarterly sheets).
$q1.Range("A2").Value = 0

$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = "515590"
$q1.Range("C2").Value = "前海开源中证500等权重ETF"
$q1.Range("D2").Value = "0.38"
$q1.Range("E2").Value = "95.07"
$q1.Range("F2").Value = "0.33"
$q1.Range("G2").Value = "0.0013"
$q1.Range("H2").Value = 9

# Formatting matching the other quarterly sheets: bold, thin border,
# centered/top-aligned header row and index column.
$q1Fmt = $q1.Range("B1:H1")
$q1Fmt.Font.Bold = $true
$q1Fmt.Borders.LineStyle = 1
$q1Fmt.HorizontalAlignment = -4108
$q1Fmt.VerticalAlignment = -4160

$q1IndexFmt = $q1.Range("A2")
$q1IndexFmt.Font.Bold = $true
$q1IndexFmt.Borders.LineStyle = 1
$q1IndexFmt.HorizontalAlignment = -4108
$q1IndexFmt.VerticalAlignment = -4160

# ---------------------------------------------------------------------
# 2) Update "总计" sheet with the new 2022-Q1 row on top
# ---------------------------------------------------------------------
# Re-fetch by name: inserting/adding sheets shifts tab positions, and we
# need the actual "总计" sheet here (not a stale position-bound handle).
$totals = $wb.Worksheets.Item("总计")
$totals.Rows.Item(2).Insert()
$totals.Range("B2:D2").ClearFormats()

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 0

$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2

$totalsIndexFmt = $totals.Range("A2")
$totalsIndexFmt.Font.Bold = $true
$totalsIndexFmt.Borders.LineStyle = 1
$totalsIndexFmt.HorizontalAlignment = -4108
$totalsIndexFmt.VerticalAlignment = -4160
